$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Split the intro paragraph into two paragraphs:
#    "An Airbnb co-host has approached you to build a web application to
#    help them manage their bookings. A property may be rented to many
#    guests at a time. A guest is registered with only one holiday
#    property."
#    becomes two paragraphs with new wording for the first sentence.
# ---------------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(2)
$introPara.Range.Text = "A holidays property manager (co-host) has approached you to build a web application to help them manage their bookings, guests and properties on a digital agenda. "
$introPara.Range.InsertParagraphAfter()
$secondPara = $d.Paragraphs.Item(3)
$secondPara.Range.Text = "A property may be rented to many guests at a time. A guest is registered with only one holiday property."

# ---------------------------------------------------------------------------
# 2. Shorten the "register / track guests" sentence (drop the trailing
#    "Important information..." clause - the detailed list below it is
#    being removed entirely).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The co-host wants to be able to register / track guests. Important information for the co-host to know is -",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The co-host wants to be able to register / track guests", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Remove the six sub-bullets that used to itemise the guest information
#    the co-host should track: Name, Last name, Date Of Birth, Nationality,
#    Contact details, Documents. After step 1 these paragraphs are numbers
#    6-11 (originally 5-10, shifted by the extra paragraph inserted above).
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(6)
$endPara = $d.Paragraphs.Item(11)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# ---------------------------------------------------------------------------
# 4. Replace the first "Possible Extensions" bullet - it now talks about a
#    search tab with autocomplete instead of sorted views. After steps 1-3
#    this is paragraph 10.
# ---------------------------------------------------------------------------
$searchPara = $d.Paragraphs.Item(10)
$searchPara.Range.Text = "There should be a search tab with autocomplete enabling the property manager to search for guests and properties"

# ---------------------------------------------------------------------------
# 5. Insert a new bullet (with the corrected "parameters" spelling) right
#    before the "If a property has multiple guests..." bullet (paragraph 11)
#    - this is the sorted-views sentence that used to occupy the first
#    bullet slot.
# ---------------------------------------------------------------------------
$multiGuestPara = $d.Paragraphs.Item(11)
$multiGuestPara.Range.InsertParagraphBefore()
$newBulletPara = $d.Paragraphs.Item(11)
$newBulletPara.Range.Text = "The views should allow the co-host to see bookings, guests and properties sorted by all parameters"

# ---------------------------------------------------------------------------
# 6. Add a trailing period to the final "Add extra functionality..." bullet.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Add extra functionality of your choosing - calculating the earning per booking after commissions, providing services (late check-in, bike rentals etc)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Add extra functionality of your choosing - calculating the earning per booking after commissions, providing services (late check-in, bike rentals etc).", 2) | Out-Null
